# Updated on Feb 07 - append three new confirmed-case rows (31, 32, 33)
# to the Singapore COVID-19 case-tracking sheet, and move the viewport /
# active selection down to the newly added data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Give the three new rows (32:34) the same look as the existing data
# rows: row height 21, 16pt non-bold font, and column D formatted as
# text (@) so date-like strings such as "Feb-07" are stored verbatim.
# ---------------------------------------------------------------------
$newRows = $ws.Range("A32:K34")
$newRows.RowHeight = 21
$newRows.Font.Size = 16
$newRows.Font.Bold = $false
$ws.Range("D32:D34").NumberFormat = "@"

# ---------------------------------------------------------------------
# Row 32 - Case 31
# ---------------------------------------------------------------------
$ws.Range("A32").Value = 31
$ws.Range("B32").Value = 1.355457
$ws.Range("C32").Value = 103.951563
$ws.Range("D32").Value = "Feb-07"
$ws.Range("E32").Value = 53
$ws.Range("F32").Value = "Male"
$ws.Range("G32").Value = "Singapore"
$ws.Range("H32").Value = "Tampines Street 24"
$ws.Range("I32").Value = "The Life Church, Missions Singapore"

# ---------------------------------------------------------------------
# Row 33 - Case 32
# (Stay / column I is entered before Visited / column H, matching the
# original shared-string insertion order.)
# ---------------------------------------------------------------------
$ws.Range("A33").Value = 32
$ws.Range("B33").Value = 1.378105
$ws.Range("C33").Value = 103.94172
$ws.Range("D33").Value = "Feb-07"
$ws.Range("E33").Value = 42
$ws.Range("F33").Value = "Female"
$ws.Range("G33").Value = "Singapore"
$ws.Range("I33").Value = "Parkway East Hospital"
$ws.Range("H33").Value = "Elias Road"

# ---------------------------------------------------------------------
# Row 34 - Case 33
# (Again, Stay / column I before Visited / column H.)
# ---------------------------------------------------------------------
$ws.Range("A34").Value = 33
$ws.Range("B34").Value = 1.392699
$ws.Range("C34").Value = 103.894414
$ws.Range("D34").Value = "Feb-07"
$ws.Range("E34").Value = 39
$ws.Range("F34").Value = "Female"
$ws.Range("G34").Value = "Singapore"
$ws.Range("I34").Value = "Sengkang Polyclinic"
$ws.Range("H34").Value = "Sengkang"

# ---------------------------------------------------------------------
# Scroll the view down to the new rows and move the active selection.
# ---------------------------------------------------------------------
$excel.ActiveWindow.ScrollRow = 18
$excel.ActiveWindow.ScrollColumn = 1
$ws.Range("E36").Select()
